# Apply cryptos list update (Thu Oct  5 04:35:58 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.688.64"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "1.644.16"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("E10").Value = "  +0.63%  "

$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").Value = "1.877.05"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "1.628.59"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").Value = "27.665.89"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("E19").Value = "  +0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.01%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.20%  "

$ws.Range("E24").Value = "  -2.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0487"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("D33").Value = "1.442.68"
$ws.Range("E33").Value = "  +2.26%  "

$ws.Range("E34").Value = "  +1.12%  "

$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("E36").Value = "  -1.15%  "

$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.887"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.57%  "

$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("E43").Value = "  +3.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.17%  "

$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").Value = "1.786.50"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("E47").Value = "  +5.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("E51").Value = "  +0.84%  "
